$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BLP-URLs")

# Insert a new row at position 2, shifting existing rows (and the data below) down
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with data for "Juist, Inselgemeinde"
$ws.Cells.Item(2, 1).Value = 452013
$ws.Cells.Item(2, 2).Value = "Juist, Inselgemeinde"
$ws.Cells.Item(2, 3).Value = 53.678347
$ws.Cells.Item(2, 4).Value = 6.995328
$ws.Cells.Item(2, 5).Value = "https://drive.google.com/drive/folders/0BxMfdWAA8UdsUFhRMzdDdnh4Z1E"
$ws.Cells.Item(2, 6).Value = "http://www.gemeinde-juist.de/"

# Extend the _FilterDatabase defined name by one row to account for the inserted row
$fdb = $wb.Names.Item("BLP-URLs!_FilterDatabase")
$fdb.RefersTo = "='BLP-URLs'!`$A`$1:`$G`$405"

# Select the newly inserted row, matching the resulting selection state
$ws.Rows.Item(2).Select()
